$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with fresh price/volume data (and two coin-name swaps)
# Each text cell is written via NumberFormat="@" so Excel keeps it as a literal
# string instead of re-parsing look-alike numbers (e.g. "0.527", "7.32"),
# then the cell style is reset to "Normal" so no stray numFmt/style is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "59.723.85"
Set-TextValue "E2" "  +3.30%  "

Set-TextValue "D3" "3.198.32"
Set-TextValue "E3" "  +2.22%  "

Set-TextValue "D5" "535.38"
Set-TextValue "E5" "  +0.40%  "

Set-TextValue "D6" "143.37"
Set-TextValue "E6" "  +3.24%  "

Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.09%  "

Set-TextValue "D8" "0.527"
Set-TextValue "E8" "  +10.84%  "

Set-TextValue "D9" "7.32"
Set-TextValue "E9" "  -0.24%  "

Set-TextValue "D10" "0.440"
Set-TextValue "E10" "  +6.66%  "

Set-TextValue "D11" "0.113"
Set-TextValue "E11" "  +4.77%  "

Set-TextValue "D12" "3.742.65"
Set-TextValue "E12" "  +2.09%  "

Set-TextValue "E13" "  +1.79%  "

Set-TextValue "D14" "26.15"
Set-TextValue "E14" "  +0.84%  "

Set-TextValue "D15" "0.0000174"
Set-TextValue "E15" "  +5.70%  "

Set-TextValue "D16" "59.730.30"
Set-TextValue "E16" "  +3.15%  "

Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.197.80"
Set-TextValue "E17" "  +2.27%  "

Set-TextValue "B18" "Polkadot"
Set-TextValue "C18" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "6.28"
Set-TextValue "E18" "  +3.14%  "

Set-TextValue "D19" "13.11"
Set-TextValue "E19" "  +3.10%  "

Set-TextValue "D20" "8.24"
Set-TextValue "E20" "  +1.80%  "

Set-TextValue "D21" "377.26"
Set-TextValue "E21" "  +2.55%  "

Set-TextValue "D22" "1.00"
Set-TextValue "E22" "  +0.14%  "

Set-TextValue "D23" "0.535"
Set-TextValue "E23" "  +5.61%  "

Set-TextValue "D24" "70.35"
Set-TextValue "E24" "  +1.58%  "

Set-TextValue "D25" "0.169"
Set-TextValue "E25" "  +0.38%  "

Set-TextValue "E26" "  -0.45%  "

Set-TextValue "D27" "8.46"
Set-TextValue "E27" "  +15.81%  "

Set-TextValue "D28" "0.0₃0887"
Set-TextValue "E28" "  +2.51%  "

Set-TextValue "D29" "22.50"
Set-TextValue "E29" "  +4.82%  "

Set-TextValue "E30" "  +1.44%  "

Set-TextValue "D31" "6.16"
Set-TextValue "E31" "  +0.92%  "

Set-TextValue "D32" "5.33"
Set-TextValue "E32" "  +3.34%  "

Set-TextValue "D33" "1.18"
Set-TextValue "E33" "  +0.86%  "

Set-TextValue "D34" "6.40"
Set-TextValue "E34" "  +5.09%  "

Set-TextValue "D35" "157.72"
Set-TextValue "E35" "  -1.12%  "

Set-TextValue "E36" "  +4.31%  "

Set-TextValue "B37" "Hedera"
Set-TextValue "C37" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D37" "0.0717"
Set-TextValue "E37" "  +6.63%  "

Set-TextValue "B38" "EnergySwap"
Set-TextValue "C38" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D38" "25.64"
Set-TextValue "E38" "  +0.79%  "

Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.73"
Set-TextValue "E39" "  +3.01%  "

Set-TextValue "B40" "Maker"
Set-TextValue "C40" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "2.711.78"
Set-TextValue "E40" "  +7.15%  "

Set-TextValue "D41" "4.30"
Set-TextValue "E41" "  +5.18%  "

Set-TextValue "D42" "0.730"
Set-TextValue "E42" "  +4.42%  "

Set-TextValue "E43" "  +9.17%  "

Set-TextValue "D44" "39.34"
Set-TextValue "E44" "  +4.27%  "

Set-TextValue "E45" "  +0.10%  "

Set-TextValue "D46" "3.236.34"
Set-TextValue "E46" "  +2.08%  "

Set-TextValue "D47" "0.998"
Set-TextValue "E47" "  +2.21%  "

Set-TextValue "E48" "  +11.91%  "

Set-TextValue "E49" "  +1.77%  "

Set-TextValue "D50" "20.61"
Set-TextValue "E50" "  +4.19%  "

Set-TextValue "D51" "0.766"
Set-TextValue "E51" "  +2.78%  "
